$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3982.5557
$ws.Range("I40").Value = 3009.2307
$ws.Range("J40").Value = 4886.357
$ws.Range("K40").Value = 3009.2307
$ws.Range("L40").Value = 4886.357
$ws.Range("M40").Value = -2834.2307
$ws.Range("N40").Value = -5236.357

# Row 41
$ws.Range("H41").Value = 8118.7334
$ws.Range("J41").Value = 16148.714
$ws.Range("L41").Value = 16148.714
$ws.Range("N41").Value = -17028.714

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = ""

# Row 96
$ws.Range("H96").Value = 479.44446
$ws.Range("I96").Value = 658.5
$ws.Range("K96").Value = 1975.5
$ws.Range("M96").Value = -602.5

# Row 97
$ws.Range("H97").Value = 1238.9166
$ws.Range("J97").Value = 1608.375
$ws.Range("L97").Value = 4825.125
$ws.Range("N97").Value = -5817.125

# Row 111
$ws.Range("H111").Value = 3261.0527
$ws.Range("I111").Value = 2782.4614
$ws.Range("K111").Value = 8347.3842
$ws.Range("M111").Value = -5280.3842

# Row 129
$ws.Range("H129").Value = 169348.5
$ws.Range("I129").Value = 252161.38
$ws.Range("J129").Value = 3722.75
$ws.Range("K129").Value = 756484.14
$ws.Range("L129").Value = 11168.25
$ws.Range("M129").Value = -751484.14
$ws.Range("N129").Value = -21168.25

# Row 132
$ws.Range("H132").Value = 15512.046
$ws.Range("I132").Value = 3419.5134
$ws.Range("K132").Value = 10258.5402
$ws.Range("M132").Value = -7728.540199999999

# Row 137
$ws.Range("H137").Value = 3547.6985
$ws.Range("I137").Value = 3723.04
$ws.Range("J137").Value = 2873.3076
$ws.Range("K137").Value = 11169.12
$ws.Range("L137").Value = 8619.9228
$ws.Range("M137").Value = -8619.119999999999
$ws.Range("N137").Value = -13719.9228


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12350475
$ws.Range("I32").Value = 13703587
$ws.Range("J32").Value = 3325.5
$ws.Range("K32").Value = 13703587
$ws.Range("L32").Value = 3325.5
$ws.Range("M32").Value = -13703300
$ws.Range("N32").Value = -3899.5

# Row 61
$ws.Range("H61").Value = 2393.6785
$ws.Range("I61").Value = 1539.25
$ws.Range("K61").Value = 1539.25
$ws.Range("M61").Value = -1327.25

# Row 97
$ws.Range("H97").Value = 4833488
$ws.Range("I97").Value = 1890.35
$ws.Range("K97").Value = 1890.35
$ws.Range("M97").Value = -1394.35

# Row 132
$ws.Range("H132").Value = 1336.7778
$ws.Range("I132").Value = 1283.2667
$ws.Range("K132").Value = 3849.800099999999
$ws.Range("M132").Value = -1319.800099999999

# Row 136
$ws.Range("H136").Value = 2393.6785
$ws.Range("I136").Value = 1539.25
$ws.Range("K136").Value = 4617.75
$ws.Range("M136").Value = -2067.75


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1783
$ws.Range("I20").Value = 1618.25
$ws.Range("J20").Value = 2002.6666
$ws.Range("K20").Value = 1618.25
$ws.Range("L20").Value = 2002.6666
$ws.Range("M20").Value = -1371.25
$ws.Range("N20").Value = -2496.6666

# Row 95
$ws.Range("H95").Value = 1125
$ws.Range("J95").Value = 1125
$ws.Range("L95").Value = 1125
$ws.Range("N95").Value = -6617

# Row 99
$ws.Range("H99").Value = 26572.588
$ws.Range("I99").Value = 39133.547
$ws.Range("K99").Value = 39133.547
$ws.Range("M99").Value = -37635.547

# Row 105
$ws.Range("H105").Value = 3224.9092
$ws.Range("I105").Value = 2072.4443
$ws.Range("K105").Value = 2072.4443
$ws.Range("M105").Value = -325.4443000000001


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 894.4286
$ws.Range("J15").Value = 165
$ws.Range("L15").Value = 165
$ws.Range("N15").Value = -505

# Row 31
$ws.Range("H31").Value = 3156.5715
$ws.Range("I31").Value = 2532.9
$ws.Range("K31").Value = 2532.9
$ws.Range("M31").Value = -2237.9

# Row 34
$ws.Range("H34").Value = 3156.5715
$ws.Range("I34").Value = 2532.9
$ws.Range("K34").Value = 2532.9
$ws.Range("M34").Value = -2330.9

# Row 86
$ws.Range("H86").Value = 47397.3
$ws.Range("J86").Value = 30242.25
$ws.Range("L86").Value = 30242.25
$ws.Range("N86").Value = -32488.25

# Row 89
$ws.Range("H89").Value = 47397.3
$ws.Range("J89").Value = 30242.25
$ws.Range("L89").Value = 151211.25
$ws.Range("N89").Value = -162443.25

# Row 99
$ws.Range("H99").Value = 28054208
$ws.Range("I99").Value = 6098789
$ws.Range("K99").Value = 6098789
$ws.Range("M99").Value = -6097291

# Row 107
$ws.Range("H107").Value = 13049.883
$ws.Range("I107").Value = 1192.25
$ws.Range("J107").Value = 23590
$ws.Range("K107").Value = 1192.25
$ws.Range("L107").Value = 23590
$ws.Range("M107").Value = 727.75
$ws.Range("N107").Value = -27430

# Row 120
$ws.Range("H120").Value = 50324
$ws.Range("J120").Value = 50324
$ws.Range("L120").Value = 50324
$ws.Range("N120").Value = -57582

# Row 122
$ws.Range("H122").Value = 447494.22
$ws.Range("J122").Value = 8499.875
$ws.Range("L122").Value = 25499.625
$ws.Range("N122").Value = -30399.625

# Row 126
$ws.Range("H126").Value = 28054208
$ws.Range("I126").Value = 6098789
$ws.Range("K126").Value = 18296367
$ws.Range("M126").Value = -18293897

# Row 132
$ws.Range("H132").Value = 3006.111
$ws.Range("I132").Value = 3006.9375
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 9020.8125
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -6490.8125
$ws.Range("N132").Value = -14058.5

# Row 134
$ws.Range("H134").Value = 4022.5
$ws.Range("I134").Value = 4199.294
$ws.Range("J134").Value = 3791.3076
$ws.Range("K134").Value = 12597.882
$ws.Range("L134").Value = 11373.9228
$ws.Range("M134").Value = -10062.882
$ws.Range("N134").Value = -16443.9228


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("I4").Value = 45727856
$ws.Range("J4").Value = 442.16666
$ws.Range("K4").Value = 137183568
$ws.Range("L4").Value = 1326.49998
$ws.Range("M4").Value = -137183456
$ws.Range("N4").Value = -1550.49998

# Row 8
$ws.Range("H8").Value = 3365.3333
$ws.Range("I8").Value = 3365.3333
$ws.Range("K8").Value = 10095.9999
$ws.Range("M8").Value = -9956.999899999999

# Row 95
$ws.Range("H95").Value = 6666.6665
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 6666.6665
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 19999.9995
$ws.Range("M95").Value = ""
$ws.Range("N95").Value = -24117.9995

# Row 132
$ws.Range("H132").Value = 1499
$ws.Range("I132").Value = 1598.6666
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 14387.9994
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -11857.9994
$ws.Range("N132").Value = -15860


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 92
$ws.Range("H92").Value = 17487.5
$ws.Range("J92").Value = 17487.5
$ws.Range("L92").Value = 17487.5
$ws.Range("N92").Value = -21231.5

# Row 122
$ws.Range("H122").Value = 4276.5
$ws.Range("I122").Value = 3406.4707
$ws.Range("J122").Value = 5146.5293
$ws.Range("K122").Value = 10219.4121
$ws.Range("L122").Value = 15439.5879
$ws.Range("M122").Value = -7769.4121
$ws.Range("N122").Value = -20339.5879

# Row 132
$ws.Range("H132").Value = 3106.111
$ws.Range("I132").Value = 3103
$ws.Range("J132").Value = 3119.8
$ws.Range("K132").Value = 9309
$ws.Range("L132").Value = 9359.400000000001
$ws.Range("M132").Value = -6779
$ws.Range("N132").Value = -14419.4


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1151.6818
$ws.Range("I61").Value = 1063.6666
$ws.Range("K61").Value = 1063.6666
$ws.Range("M61").Value = -861.6666

# Row 113
$ws.Range("H113").Value = 1151.6818
$ws.Range("I113").Value = 1063.6666
$ws.Range("K113").Value = 1063.6666
$ws.Range("M113").Value = 1106.3334

# Row 122
$ws.Range("H122").Value = 7366.2666
$ws.Range("I122").Value = 5622
$ws.Range("J122").Value = 7634.615
$ws.Range("K122").Value = 16866
$ws.Range("L122").Value = 22903.845
$ws.Range("M122").Value = -14416
$ws.Range("N122").Value = -27803.845

# Row 125
$ws.Range("H125").Value = 67499.914
$ws.Range("J125").Value = 67499.914
$ws.Range("L125").Value = 67499.914
$ws.Range("N125").Value = -77339.914

# Row 132
$ws.Range("H132").Value = 6105.923
$ws.Range("I132").Value = 3873.6667
$ws.Range("J132").Value = 9149.909
$ws.Range("K132").Value = 11621.0001
$ws.Range("L132").Value = 27449.727
$ws.Range("M132").Value = -9091.000100000001
$ws.Range("N132").Value = -32509.727

# Row 136
$ws.Range("H136").Value = 2946.6853
$ws.Range("I136").Value = 2652.551
$ws.Range("K136").Value = 7957.653
$ws.Range("M136").Value = -5407.653


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 579.3461
$ws.Range("I107").Value = 652.1
$ws.Range("J107").Value = 336.83334
$ws.Range("K107").Value = 1956.3
$ws.Range("L107").Value = 1010.50002
$ws.Range("M107").Value = -36.30000000000018
$ws.Range("N107").Value = -4850.50002

# Row 122
$ws.Range("H122").Value = 269889.75
$ws.Range("I122").Value = 2397.5483
$ws.Range("J122").Value = 1023731.44
$ws.Range("K122").Value = 7192.644899999999
$ws.Range("L122").Value = 3071194.32
$ws.Range("M122").Value = -4742.644899999999
$ws.Range("N122").Value = -3076094.32

# Row 126
$ws.Range("H126").Value = 2799.25
$ws.Range("I126").Value = 1218.6
$ws.Range("K126").Value = 3655.8
$ws.Range("M126").Value = -1185.8

# Row 136
$ws.Range("H136").Value = 2351.5386
$ws.Range("I136").Value = 1355.7142
$ws.Range("K136").Value = 4067.1426
$ws.Range("M136").Value = -1517.1426

